$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "612×3=1836"
$t.Cell(1, 2).Range.Text = "528×7=3696"
$t.Cell(1, 3).Range.Text = "653×2=1306"
$t.Cell(1, 4).Range.Text = "469×3=1407"
$t.Cell(1, 5).Range.Text = "143×5=715"

$t.Cell(5, 1).Range.Text = "503×4=2012"
$t.Cell(5, 2).Range.Text = "562×9=5058"
$t.Cell(5, 3).Range.Text = "376×9=3384"
$t.Cell(5, 4).Range.Text = "774×4=3096"
$t.Cell(5, 5).Range.Text = "982×5=4910"

$t.Cell(10, 1).Range.Text = "483×5=2415"
$t.Cell(10, 2).Range.Text = "110×3=330"
$t.Cell(10, 3).Range.Text = "966×6=5796"
$t.Cell(10, 4).Range.Text = "113×9=1017"
$t.Cell(10, 5).Range.Text = "660×8=5280"

$t.Cell(15, 1).Range.Text = "217×6=1302"
$t.Cell(15, 2).Range.Text = "697×2=1394"
$t.Cell(15, 3).Range.Text = "848×7=5936"
$t.Cell(15, 4).Range.Text = "538×2=1076"
$t.Cell(15, 5).Range.Text = "848×3=2544"

$t.Cell(20, 1).Range.Text = "513×4=2052"
$t.Cell(20, 2).Range.Text = "700×7=4900"
$t.Cell(20, 3).Range.Text = "945×4=3780"
$t.Cell(20, 4).Range.Text = "499×8=3992"
$t.Cell(20, 5).Range.Text = "905×9=8145"
